# Rename the worksheet tabs (commit: "renamed tabs of dn80 bacterial go terms")
$wb = $excel.ActiveWorkbook

$wsRaw      = $wb.Worksheets.Item(1)   # was: DN80_combined_TD+ND_bp_GO
$wsRevigo   = $wb.Worksheets.Item(2)   # was: revigo
$wsOrganize = $wb.Worksheets.Item(3)   # was: organizing
$wsPython   = $wb.Worksheets.Item(4)   # was: for python
$wsGoTerms  = $wb.Worksheets.Item(5)   # GO terms >1 peptides (name unchanged)

$wsRaw.Name      = "Raw GO counts DN80_combined_TD+"
$wsRevigo.Name   = "Revigo condensation"
$wsOrganize.Name = "Manual condensation after Revig"
$wsPython.Name   = "condensed for Fig. 6b python sc"

# Update the header labels on the "condensed for Fig. 6b python sc" sheet
$wsPython.Range("B1").Value = "Day 0 Bacteria"
$wsPython.Range("C1").Value = "Day 2 Bacteria"
$wsPython.Range("D1").Value = "Day 5 Bacteria"
$wsPython.Range("E1").Value = "Day 12 Bacteria"

# Restore that sheet's own view/selection to the top (A1 -> B1), no longer the active tab
$wsPython.Activate()
$wsPython.Range("B1").Select()

# The active tab moves from "condensed for Fig. 6b python sc" (index 4) to "GO terms >1 peptides" (index 5)
$wsGoTerms.Activate()
$wsGoTerms.Range("C226").Select()
